$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a string value to a cell while guaranteeing it is stored as
# text (Excel would otherwise auto-coerce numeric-looking strings like
# "1.00" or "20.00" into numbers). We temporarily force a Text number
# format, assign the value, then restore the cell to the default/
# "Normal" style so no visible formatting/style change is left behind.
function Set-TextValue($ref, $text) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "51.256.40"
Set-TextValue "E2" "  -1.31%  "
Set-TextValue "D3" "2.770.19"
Set-TextValue "E3" "  -0.15%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "353.39"
Set-TextValue "E5" "  -0.68%  "
Set-TextValue "D6" "107.53"
Set-TextValue "E6" "  -1.21%  "
Set-TextValue "E7" "  -2.65%  "
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.582"
Set-TextValue "E9" "  -1.16%  "
Set-TextValue "D10" "39.46"
Set-TextValue "E10" "  -1.70%  "
Set-TextValue "E11" "  +3.34%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D12" "20.00"
Set-TextValue "E12" "  +3.35%  "
$ws.Range("B13").Value = "Dogecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D13" "0.0830"
Set-TextValue "E13" "  -2.49%  "
Set-TextValue "E14" "  -1.01%  "
Set-TextValue "D15" "3.204.62"
Set-TextValue "E15" "  -0.16%  "
Set-TextValue "D16" "2.768.49"
Set-TextValue "E16" "  +0.06%  "
Set-TextValue "E17" "  -0.75%  "
Set-TextValue "D18" "51.199.27"
Set-TextValue "E18" "  -1.13%  "
Set-TextValue "D19" "7.64"
Set-TextValue "E19" "  +3.24%  "
Set-TextValue "D20" "3.11"
Set-TextValue "E20" "  -0.52%  "
Set-TextValue "D21" "13.08"
Set-TextValue "E21" "  +0.50%  "
Set-TextValue "D22" "0.0₃0959"
Set-TextValue "E22" "  -1.58%  "
Set-TextValue "D23" "69.58"
Set-TextValue "E23" "  -0.31%  "
Set-TextValue "D24" "265.28"
Set-TextValue "E24" "  -3.26%  "
Set-TextValue "E25" "  -0.70%  "
Set-TextValue "E26" "  -0.04%  "
Set-TextValue "E27" "  -2.21%  "
Set-TextValue "E28" "  +13.13%  "
Set-TextValue "D29" "10.15"
Set-TextValue "E29" "  +0.29%  "
Set-TextValue "D30" "2.20"
Set-TextValue "E30" "  -0.67%  "
Set-TextValue "D31" "35.73"
Set-TextValue "E31" "  +6.01%  "
Set-TextValue "D32" "51.80"
Set-TextValue "E32" "  +0.43%  "
Set-TextValue "E33" "  +6.55%  "
Set-TextValue "D34" "5.55"
Set-TextValue "E34" "  +4.20%  "
Set-TextValue "E35" "  -4.76%  "
Set-TextValue "D36" "0.0826"
Set-TextValue "E36" "  -1.86%  "
Set-TextValue "E37" "  +0.07%  "
Set-TextValue "D38" "18.13"
Set-TextValue "E38" "  +0.64%  "
Set-TextValue "E39" "  -2.19%  "
Set-TextValue "E40" "  -1.64%  "
Set-TextValue "E41" "  -0.31%  "
Set-TextValue "E42" "  -1.20%  "
Set-TextValue "D43" "121.14"
Set-TextValue "E43" "  -0.20%  "
Set-TextValue "D44" "22.09"
Set-TextValue "E44" "  +0.38%  "
Set-TextValue "E45" "  -2.39%  "
Set-TextValue "D46" "2.095.30"
Set-TextValue "E46" "  +2.14%  "
Set-TextValue "D47" "3.23"
Set-TextValue "E47" "  -0.39%  "
Set-TextValue "D48" "2.29"
Set-TextValue "E48" "  +1.06%  "
Set-TextValue "D49" "0.905"
Set-TextValue "E49" "  -1.60%  "
Set-TextValue "E50" "  -4.95%  "
Set-TextValue "E51" "  +6.88%  "
